$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to match the player
$ws.Name = "Trent Boult"

# Force the whole used range to be stored as text so that
# numeric-looking values ("0", "1", "100.00", ...) are kept
# as literal text, same as the rest of the scraped columns.
$ws.Range("A1:M5").NumberFormat = "@"

# Header row (a new "matchNo" column was inserted before "teamName")
$ws.Cells.Item(1, 1).Value = "matchNo"
$ws.Cells.Item(1, 2).Value = "teamName"
$ws.Cells.Item(1, 3).Value = "batterName"
$ws.Cells.Item(1, 4).Value = "states"
$ws.Cells.Item(1, 5).Value = "runs"
$ws.Cells.Item(1, 6).Value = "balls"
$ws.Cells.Item(1, 7).Value = "fours"
$ws.Cells.Item(1, 8).Value = "sixes"
$ws.Cells.Item(1, 9).Value = "sr"
$ws.Cells.Item(1, 10).Value = "opponentTeamName"
$ws.Cells.Item(1, 11).Value = "venue"
$ws.Cells.Item(1, 12).Value = "date"
$ws.Cells.Item(1, 13).Value = "result"

# Data rows: one row per innings scraped for Trent Boult
# Row 2
$ws.Cells.Item(2, 1).Value = "55th"
$ws.Cells.Item(2, 2).Value = "Mumbai Indians"
$ws.Cells.Item(2, 3).Value = "Trent Boult"
$ws.Cells.Item(2, 4).Value = ""
$ws.Cells.Item(2, 5).Value = "0"
$ws.Cells.Item(2, 6).Value = "0"
$ws.Cells.Item(2, 7).Value = "0"
$ws.Cells.Item(2, 8).Value = "0"
$ws.Cells.Item(2, 9).Value = "-"
$ws.Cells.Item(2, 10).Value = "Sunrisers Hyderabad"
$ws.Cells.Item(2, 11).Value = "Abu Dhabi"
$ws.Cells.Item(2, 12).Value = "October 08"
$ws.Cells.Item(2, 13).Value = "Mumbai won by 42 runs"

# Row 3
$ws.Cells.Item(3, 1).Value = "13th"
$ws.Cells.Item(3, 2).Value = "Mumbai Indians"
$ws.Cells.Item(3, 3).Value = "Trent Boult"
$ws.Cells.Item(3, 4).Value = ""
$ws.Cells.Item(3, 5).Value = "1"
$ws.Cells.Item(3, 6).Value = "1"
$ws.Cells.Item(3, 7).Value = "0"
$ws.Cells.Item(3, 8).Value = "0"
$ws.Cells.Item(3, 9).Value = "100.00"
$ws.Cells.Item(3, 10).Value = "Delhi Capitals"
$ws.Cells.Item(3, 11).Value = "Chennai"
$ws.Cells.Item(3, 12).Value = "April 20"
$ws.Cells.Item(3, 13).Value = "Capitals won by 6 wickets (with 5 balls remaining)"

# Row 4
$ws.Cells.Item(4, 1).Value = "5th"
$ws.Cells.Item(4, 2).Value = "Mumbai Indians"
$ws.Cells.Item(4, 3).Value = "Trent Boult"
$ws.Cells.Item(4, 4).Value = ""
$ws.Cells.Item(4, 5).Value = "0"
$ws.Cells.Item(4, 6).Value = "0"
$ws.Cells.Item(4, 7).Value = "0"
$ws.Cells.Item(4, 8).Value = "0"
$ws.Cells.Item(4, 9).Value = "-"
$ws.Cells.Item(4, 10).Value = "Kolkata Knight Riders"
$ws.Cells.Item(4, 11).Value = "Chennai"
$ws.Cells.Item(4, 12).Value = "April 13"
$ws.Cells.Item(4, 13).Value = "Mumbai won by 10 runs"

# Row 5
$ws.Cells.Item(5, 1).Value = "39th"
$ws.Cells.Item(5, 2).Value = "Mumbai Indians"
$ws.Cells.Item(5, 3).Value = "Trent Boult"
$ws.Cells.Item(5, 4).Value = ""
$ws.Cells.Item(5, 5).Value = "0"
$ws.Cells.Item(5, 6).Value = "3"
$ws.Cells.Item(5, 7).Value = "0"
$ws.Cells.Item(5, 8).Value = "0"
$ws.Cells.Item(5, 9).Value = "0.00"
$ws.Cells.Item(5, 10).Value = "Royal Challengers Bangalore"
$ws.Cells.Item(5, 11).Value = "Dubai (DSC)"
$ws.Cells.Item(5, 12).Value = "September 26"
$ws.Cells.Item(5, 13).Value = "RCB won by 54 runs"

